# Automatische test-sync: 2025-06-22 22:07:50
# Append the new "BTW-nummer toevoegen" log entry to the Logs sheet (row 55),
# extend the conditional formatting ranges to cover the new row, and update
# the Dashboard category-count table to reflect the new "Factuur / Administratie"
# entry (which now has 3 occurrences and moves above "Bestelling / Levering").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: add new row 55
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(55, 1).Value = "BTW-nummer toevoegen"
$logs.Cells.Item(55, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(55, 3).Value = "Mijn BTW-nummer is niet vermeld op de factuur."
$logs.Cells.Item(55, 4).Value = "Factuur / Administratie"
$logs.Cells.Item(55, 5).Value = "Geachte klant,`nDank u voor uw bericht. Om uw probleem met betrekking tot het ontbrekende BTW-nummer op de factuur op te lossen, ontvangen wij graag meer informatie om uw specifieke situatie te begrijpen. Kunt u ons uw factuurnummer en bedrijfsnaam doorgeven, zodat wij dit verder kunnen onderzoeken en indien nodig corrigeren?`nWij streven ernaar om u zo snel mogelijk van dienst te zijn en verontschuldigen ons voor het ongemak.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item(55, 6).Value = "2025-06-22 22:07:17"
$logs.Cells.Item(55, 7).Value = "Ja"

# Re-fit the row height so it matches the sheet's default (no explicit/custom
# row height), the same way the other data rows are stored.
$logs.Rows.Item(55).AutoFit()

# ---------------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges to include row 55
# ---------------------------------------------------------------------------
$dRules = $logs.Range("D2:D54").FormatConditions
for ($i = 1; $i -le $dRules.Count; $i++) {
  $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D55"))
}

$gRules = $logs.Range("G2:G54").FormatConditions
for ($i = 1; $i -le $gRules.Count; $i++) {
  $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G55"))
}

# ---------------------------------------------------------------------------
# 3. Dashboard sheet: update category counts (rows 9-12) so the table stays
#    sorted by count descending now that "Factuur / Administratie" has 3
#    occurrences.
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(9, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(9, 2).Value = 3

$dash.Cells.Item(10, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(10, 2).Value = 3

$dash.Cells.Item(11, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(11, 2).Value = 2

$dash.Cells.Item(12, 1).Value = "Overig"
$dash.Cells.Item(12, 2).Value = 2
